# Femacal de La Calera - Ajo: insert a new weekly price record before the
# existing row 572, shifting all subsequent rows down by one (old row 572
# becomes new row 573, ..., old row 640 becomes new row 641).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 572; everything below shifts down.
$ws.Rows.Item(572).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(572, 1).Value = 3
$ws.Cells.Item(572, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(572, 3).Value = "Coquimbo"
$ws.Cells.Item(572, 4).Value = 44946
$ws.Cells.Item(572, 5).Value = 5
$ws.Cells.Item(572, 6).Value = 100112003
$ws.Cells.Item(572, 7).Value = "Ajo"
$ws.Cells.Item(572, 8).Value = "Chino"
$ws.Cells.Item(572, 9).Value = "Primera"
$ws.Cells.Item(572, 10).Value = 80
$ws.Cells.Item(572, 11).Value = 14000
$ws.Cells.Item(572, 12).Value = 14500
$ws.Cells.Item(572, 13).Value = 14250
$ws.Cells.Item(572, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(572, 15).Value = "China"
$ws.Cells.Item(572, 16).Value = 1425
$ws.Cells.Item(572, 17).Value = 10
$ws.Cells.Item(572, 18).Value = "Hortaliza"

Write-Output "Inserted new row 572 for Femacal de La Calera - Ajo"
